# Update cryptos list figures (prices / 1h volume %) as scraped on
# Sun Aug 20 19:06:27 UTC 2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "Price" cell (column D) while forcing it to stay plain
# text, since many of the values look numeric (e.g. "1.011", "0.06454")
# and Excel would otherwise silently convert them to numbers / scientific
# notation. We flip the cell to Text format just long enough to assign the
# value, then put the formatting back the way it was (General / Normal
# style) so no stray style is left behind.
function Set-PriceText($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText "D2" "26.387.59"
$ws.Range("E2").Value = "  +0.32%  "

# Row 3 - Ethereum
Set-PriceText "D3" "1.692.76"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4 - TetherUSD
Set-PriceText "D4" "1.011"
$ws.Range("E4").Value = "  +0.39%  "

# Row 5 - BNB
Set-PriceText "D5" "219.14"
$ws.Range("E5").Value = "  -0.01%  "

# Row 6 - XRP
Set-PriceText "D6" "0.5482"
$ws.Range("E6").Value = "  +4.16%  "

# Row 7 - USDC
Set-PriceText "D7" "1.011"
$ws.Range("E7").Value = "  +0.34%  "

# Row 8 - Cardano
Set-PriceText "D8" "0.2736"
$ws.Range("E8").Value = "  +1.15%  "

# Row 9 - Dogecoin
Set-PriceText "D9" "0.06454"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10 - Solana
Set-PriceText "D10" "21.97"
$ws.Range("E10").Value = "  -0.32%  "

# Row 11 - TRON
Set-PriceText "D11" "0.07668"
$ws.Range("E11").Value = "  +2.62%  "

# Row 12 & 13 - Polkadot / WrappedEther swapped order
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-PriceText "D12" "1.712.39"
$ws.Range("E12").Value = "  +1.15%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-PriceText "D13" "4.541"
$ws.Range("E13").Value = "  -0.46%  "

# Row 14 - Polygon
Set-PriceText "D14" "0.5829"
$ws.Range("E14").Value = "  -0.41%  "

# Row 15 - ShibaInu
Set-PriceText "D15" "0.000008345"
$ws.Range("E15").Value = "  -2.24%  "

# Row 16 - Litecoin
Set-PriceText "D16" "65.37"
$ws.Range("E16").Value = "  +1.31%  "

# Row 17 - WrappedBTC
Set-PriceText "D17" "26.440.96"
$ws.Range("E17").Value = "  +0.30%  "

# Row 18 - Uniswap
Set-PriceText "D18" "4.937"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19 - Dai
Set-PriceText "D19" "1.011"
$ws.Range("E19").Value = "  +0.37%  "

# Row 20 - Avalanche
Set-PriceText "D20" "10.96"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21 - BitcoinCash
Set-PriceText "D21" "191.80"
$ws.Range("E21").Value = "  +1.12%  "

# Row 22 - Chainlink
Set-PriceText "D22" "6.246"
$ws.Range("E22").Value = "  +0.34%  "

# Row 23 - BinanceUSD (only E changed)
$ws.Range("E23").Value = "  +0.38%  "

# Row 24 - Monero
Set-PriceText "D24" "149.07"
$ws.Range("E24").Value = "  +2.94%  "

# Row 25 - Stellar
Set-PriceText "D25" "0.1323"
$ws.Range("E25").Value = "  +7.19%  "

# Row 26 - Cosmos
Set-PriceText "D26" "7.910"
$ws.Range("E26").Value = "  +3.10%  "

# Row 27 - EthereumClassic
Set-PriceText "D27" "15.76"
$ws.Range("E27").Value = "  -0.74%  "

# Row 28 - Hedera
Set-PriceText "D28" "0.06287"
$ws.Range("E28").Value = "  -5.66%  "

# Row 29 - Toncoin
Set-PriceText "D29" "1.385"
$ws.Range("E29").Value = "  +2.36%  "

# Row 30 - PancakeSwap (only E changed)
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - Filecoin
Set-PriceText "D31" "3.597"
$ws.Range("E31").Value = "  +0.25%  "

# Row 32 - InternetComputer(DFINITY)
Set-PriceText "D32" "3.603"
$ws.Range("E32").Value = "  +0.75%  "

# Row 33 - LidoDAOToken (only D changed)
Set-PriceText "D33" "1.684"

# Row 34 - ARBITRUM
Set-PriceText "D34" "1.043"
$ws.Range("E34").Value = "  +1.36%  "

# Row 35 - ImmutableX
Set-PriceText "D35" "0.6137"
$ws.Range("E35").Value = "  -1.41%  "

# Row 36 - HuobiToken
Set-PriceText "D36" "2.413"
$ws.Range("E36").Value = "  +0.77%  "

# Row 37 - MXToken
Set-PriceText "D37" "2.705"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38 - FraxShare
Set-PriceText "D38" "6.203"
$ws.Range("E38").Value = "  -2.44%  "

# Row 39 - VeChain
Set-PriceText "D39" "0.01639"
$ws.Range("E39").Value = "  +0.84%  "

# Row 40 - Maker
Set-PriceText "D40" "1.116.03"
$ws.Range("E40").Value = "  +0.62%  "

# Row 41 - TrustWalletToken
Set-PriceText "D41" "0.8894"
$ws.Range("E41").Value = "  +0.64%  "

# Row 42 - PaxDollar (only E changed)
$ws.Range("E42").Value = "  -0.10%  "

# Row 43 - Quant (only E changed)
$ws.Range("E43").Value = "  +1.03%  "

# Row 44 - RocketPoolETH
Set-PriceText "D44" "1.844.81"
$ws.Range("E44").Value = "  +0.33%  "

# Row 45 - BabyDogeCoin (only E changed)
$ws.Range("E45").Value = "  -3.39%  "

# Row 46 - Aave
Set-PriceText "D46" "57.51"
$ws.Range("E46").Value = "  +0.95%  "

# Row 47 & 48 - EnergySwap / Frax swapped order
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-PriceText "D47" "1.010"
$ws.Range("E47").Value = "  -0.14%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-PriceText "D48" "8.180"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49 - Cronos (only E changed)
$ws.Range("E49").Value = "  +0.29%  "

# Row 50 - Mantle
Set-PriceText "D50" "0.4305"
$ws.Range("E50").Value = "  +0.05%  "

# Row 51 - Aptos
Set-PriceText "D51" "6.092"
$ws.Range("E51").Value = "  +0.50%  "
